# Burn Down Chart: extend the tracked date range from 26 rows (A2:A27) to
# 41 rows (A2:A42). The "Hours Left" burndown (col B) keeps decrementing by
# 10.8 each day (same slope as before) but now starts higher (432 instead of
# 270) so it still reaches 0 on the new, later final day. The "Burn-Down"
# scope line (col C) simply continues flat at 69.5 for the newly added days.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Make room: insert 15 new rows after the current last data row (27) ---
$ws.Rows("28:42").Insert()

# --- 2. Column A: dates carry on day-by-day through 42102 (Nov 2015) ---
for ($r = 28; $r -le 42; $r++) {
    $ws.Cells.Item($r, 1).Value2 = 42062 + ($r - 2)
    $ws.Cells.Item($r, 1).NumberFormat = "mmmm\ d"
}

# --- 3. Column B ("Hours Left"): rebuild the whole shared-formula chain ---
#     B2:B40 = "=B{next row}+10.8" (shared formula), B41 = "=B42+10.8"
#     (the lone non-shared seed, same pattern as the old B26), B42 = literal 0.
for ($r = 2; $r -le 40; $r++) {
    $nextRow = $r + 1
    $ws.Cells.Item($r, 2).Formula = "=B$nextRow+10.8"
}
$ws.Cells.Item(41, 2).Formula = "=B42+10.8"
$ws.Cells.Item(42, 2).Value2 = 0

# The whole of column B also loses its explicit "0.0" number style as part
# of this edit, reverting to the default (General) "Normal" cell style.
$ws.Range("B2:B42").Style = "Normal"

# --- 4. Column C ("Burn-Down"): new rows just keep repeating the prior value ---
for ($r = 28; $r -le 42; $r++) {
    $prevRow = $r - 1
    $ws.Cells.Item($r, 3).Formula = "=C$prevRow"
    $ws.Cells.Item($r, 3).NumberFormat = "0.0"
}

$excel.Calculate()

# --- 5. Update the chart's two series so they reference the new A2:A42 /
#        B2:B42 / C2:C42 ranges ---
$chart = $ws.ChartObjects().Item(1).Chart
$series1 = $chart.SeriesCollection().Item(1)
$series1.Formula = '=SERIES(Sheet1!$B$1,Sheet1!$A$2:$A$42,Sheet1!$B$2:$B$42,1)'
$series2 = $chart.SeriesCollection().Item(2)
$series2.Formula = '=SERIES(Sheet1!$C$1,Sheet1!$A$2:$A$42,Sheet1!$C$2:$C$42,2)'

# The value axis number format tracked column B's old "0.0" style; now that
# B is back to General, the (source-linked) axis format follows suit.
$valueAxis = $chart.Axes(2)
$valueAxis.TickLabels.NumberFormat = "General"

$excel.Calculate()

# --- 6. Selection cosmetic change noted in the saved view state ---
$ws.Range("S8").Select()
